# Updated: po 08. 11. 2021
# Refresh AgTests (F) / AgPosit (G) figures for existing rows and append
# four new daily rows (610-613) covering 2021-11-04 .. 2021-11-07.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F503").Value = 7879
$ws.Range("F533").Value = 11965
$ws.Range("F535").Value = 10260
$ws.Range("F536").Value = 8193
$ws.Range("F537").Value = 13733
$ws.Range("F538").Value = 11349
$ws.Range("F539").Value = 10722
$ws.Range("F540").Value = 12534
$ws.Range("F541").Value = 16780
$ws.Range("F542").Value = 10430
$ws.Range("F543").Value = 4772
$ws.Range("F544").Value = 14359
$ws.Range("F545").Value = 16714
$ws.Range("F546").Value = 3972
$ws.Range("F547").Value = 14018
$ws.Range("F548").Value = 17280
$ws.Range("F549").Value = 10895
$ws.Range("F550").Value = 8581
$ws.Range("F551").Value = 17904
$ws.Range("F552").Value = 15644
$ws.Range("F553").Value = 15546
$ws.Range("G553").Value = 179
$ws.Range("F554").Value = 17990
$ws.Range("F555").Value = 21648
$ws.Range("F556").Value = 12255
$ws.Range("F557").Value = 10979
$ws.Range("F558").Value = 24722
$ws.Range("F559").Value = 22508
$ws.Range("F560").Value = 6110
$ws.Range("F561").Value = 24345
$ws.Range("F562").Value = 27117
$ws.Range("F563").Value = 14108
$ws.Range("F564").Value = 14317
$ws.Range("F565").Value = 29007
$ws.Range("F566").Value = 25758
$ws.Range("F567").Value = 23507
$ws.Range("F568").Value = 23975
$ws.Range("F569").Value = 32471
$ws.Range("F570").Value = 15229
$ws.Range("F571").Value = 15009
$ws.Range("F572").Value = 33391
$ws.Range("G572").Value = 597
$ws.Range("F573").Value = 27009
$ws.Range("F574").Value = 23448
$ws.Range("G574").Value = 356
$ws.Range("F575").Value = 26102
$ws.Range("G575").Value = 395
$ws.Range("F576").Value = 28995
$ws.Range("F577").Value = 14768
$ws.Range("F578").Value = 15086
$ws.Range("F579").Value = 32651
$ws.Range("G579").Value = 644
$ws.Range("F580").Value = 28877
$ws.Range("G580").Value = 517
$ws.Range("F581").Value = 27118
$ws.Range("G581").Value = 481
$ws.Range("F582").Value = 25963
$ws.Range("F583").Value = 29384
$ws.Range("F584").Value = 13256
$ws.Range("F585").Value = 14929
$ws.Range("G585").Value = 359
$ws.Range("F586").Value = 33672
$ws.Range("F587").Value = 28227
$ws.Range("G587").Value = 555
$ws.Range("F588").Value = 25384
$ws.Range("F589").Value = 25486
$ws.Range("G589").Value = 474
$ws.Range("F590").Value = 29337
$ws.Range("F591").Value = 14806
$ws.Range("G591").Value = 436
$ws.Range("F592").Value = 17965
$ws.Range("F593").Value = 36998
$ws.Range("G593").Value = 1194
$ws.Range("F594").Value = 29783
$ws.Range("G594").Value = 817
$ws.Range("F595").Value = 27267
$ws.Range("F596").Value = 29179
$ws.Range("F597").Value = 29514
$ws.Range("F598").Value = 15435
$ws.Range("F599").Value = 16544
$ws.Range("G599").Value = 872
$ws.Range("F600").Value = 39944
$ws.Range("G600").Value = 1674
$ws.Range("F601").Value = 31649
$ws.Range("G601").Value = 1331
$ws.Range("F602").Value = 30023
$ws.Range("G602").Value = 1289
$ws.Range("F603").Value = 31816
$ws.Range("G603").Value = 1523
$ws.Range("F604").Value = 29871
$ws.Range("G604").Value = 1526
$ws.Range("F605").Value = 14639
$ws.Range("G605").Value = 1016
$ws.Range("F606").Value = 14032
$ws.Range("G606").Value = 1243
$ws.Range("F607").Value = 10741
$ws.Range("G607").Value = 944
$ws.Range("F608").Value = 45104
$ws.Range("G608").Value = 2858
$ws.Range("F609").Value = 35826
$ws.Range("G609").Value = 2126

# New rows
$ws.Range("A610").Value = 44504
$ws.Range("B610").Value = 506795
$ws.Range("C610").Value = 20519
$ws.Range("D610").Value = 6805
$ws.Range("E610").Value = 13166
$ws.Range("F610").Value = 32812
$ws.Range("G610").Value = 1877
$ws.Range("A611").Value = 44505
$ws.Range("B611").Value = 513297
$ws.Range("C611").Value = 21305
$ws.Range("D611").Value = 6502
$ws.Range("E611").Value = 13205
$ws.Range("F611").Value = 31094
$ws.Range("G611").Value = 1919
$ws.Range("A612").Value = 44506
$ws.Range("B612").Value = 518635
$ws.Range("C612").Value = 16826
$ws.Range("D612").Value = 5338
$ws.Range("E612").Value = 13229
$ws.Range("F612").Value = 13809
$ws.Range("G612").Value = 1242
$ws.Range("A613").Value = 44507
$ws.Range("B613").Value = 521650
$ws.Range("C613").Value = 8986
$ws.Range("D613").Value = 3015
$ws.Range("E613").Value = 13269
$ws.Range("F613").Value = 16103
$ws.Range("G613").Value = 1472